$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '62.981.71'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +3.16%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.452.03'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +2.04%  '
$ws.Range('E4').Value = '  -0.21%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '577.70'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.82%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '146.17'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +2.83%  '
$ws.Range('E7').Value = '  +0.07%  '
$ws.Range('E8').Value = '  +0.18%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.452.05'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.71%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.112'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +3.08%  '
$ws.Range('E11').Value = '  +2.53%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '5.29'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +1.78%  '
$ws.Range('E13').Value = '  +2.29%  '
$ws.Range('E14').Value = '  +7.11%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.0000181'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +5.74%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.894.58'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +1.89%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '62.766.59'
$ws.Range('D17').Style = 'Normal'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.462.41'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.94%  '
$ws.Range('E19').Value = '  -1.63%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.07'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +3.48%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '333.14'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +2.56%  '
$ws.Range('E22').Value = '  +1.09%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '2.07'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +7.49%  '
$ws.Range('E24').Value = '  +0.00%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '66.44'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.83%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '652.77'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +10.62%  '
$ws.Range('E27').Value = '  +17.60%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.54'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +3.43%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.0000100'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +5.46%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.573.28'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +2.23%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '8.23'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +2.46%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.45'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +6.36%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0₆0448'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +57.79%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.88'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +3.43%  '
$ws.Range('E35').Value = '  +3.28%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.50'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +1.49%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.999'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.12%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.77'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +3.18%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.56'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +6.12%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.375'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.57%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '153.05'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.16%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '18.83'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +2.58%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.74'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +9.94%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.77'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +4.88%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '42.51'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +1.52%  '
$ws.Range('E46').Value = '  +0.02%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '15.01'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +27.57%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '145.78'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +2.68%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '3.65'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +3.28%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '20.71'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +5.23%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.606'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +2.50%  '
